$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Fix the "SPOReport" -> "SPO Reports" text on the two sheets that still
#    had the old (typo'd) label.
# ---------------------------------------------------------------------------
$devAppPools = $wb.Worksheets.Item("Dev App Pools")
$devAppPools.Range("A5").Value = "SPO Reports"

$webApps = $wb.Worksheets.Item("Web Apps")
$webApps.Range("C11").Value = "SPO Reports"

# ---------------------------------------------------------------------------
# 2) Update the selection/cursor position remembered on the "Web Apps" sheet.
# ---------------------------------------------------------------------------
$webApps.Activate()
$webApps.Range("A11").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) Add the new "Virtual Directories" sheet after "Legacy Apps" with the
#    input-definition columns/rows for creating virtual directories.
# ---------------------------------------------------------------------------
$legacyApps = $wb.Worksheets.Item("Legacy Apps")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $legacyApps)
$newSheet.Name = "Virtual Directories"

# Column widths to match the other definition sheets.
$newSheet.Columns.Item(1).ColumnWidth = 21.7109375
$newSheet.Columns.Item(2).ColumnWidth = 25.140625
$newSheet.Columns.Item(3).ColumnWidth = 35.42578125

# Header row values.
$newSheet.Range("A1").Value = "IISAppName"
$newSheet.Range("B1").Value = "vName"
$newSheet.Range("C1").Value = "DNSName"

# Data row values.
$newSheet.Range("A2").Value = "IRMA Client"
$newSheet.Range("B2").Value = "sporefiles"
$newSheet.Range("C2").Value = "E:\WebTools\PO\"

# Reuse the same black-fill / white-centered-text header formatting used by
# the other sheets' header rows (copy format only from an existing header).
$webApps.Range("A1:C1").Copy() | Out-Null
$newSheet.Range("A1:C1").PasteSpecial(-4122) | Out-Null

# Match zoom level used by the other definition sheets.
$newSheet.Activate()
$excel.ActiveWindow.Zoom = 175

# Selection left on the data row, matching the other freshly-populated sheets.
$newSheet.Range("A2:C2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4) Restore "Dev App Pools" as the active tab (unchanged from the original).
# ---------------------------------------------------------------------------
$devAppPools.Activate()
